# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) for the affected leve rows across sheets, per the latest price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 13890271
$ws.Range("I100").Value = 23810764
$ws.Range("J100").Value = 1579.2
$ws.Range("K100").Value = 23810764
$ws.Range("L100").Value = 1579.2
$ws.Range("M100").Value = -23810223
$ws.Range("N100").Value = -2661.2

$ws.Range("H112").Value = 27779746
$ws.Range("J112").Value = 35716468
$ws.Range("L112").Value = 107149404
$ws.Range("N112").Value = -107151620

$ws.Range("H129").Value = 700.3077
$ws.Range("I129").Value = 420.35
$ws.Range("J129").Value = 995
$ws.Range("K129").Value = 1261.05
$ws.Range("L129").Value = 2985
$ws.Range("M129").Value = 3738.95
$ws.Range("N129").Value = -12985

$ws.Range("H137").Value = 924.4697
$ws.Range("I137").Value = 940.38464
$ws.Range("J137").Value = 865.3570999999999
$ws.Range("K137").Value = 2821.15392
$ws.Range("L137").Value = 2596.0713
$ws.Range("M137").Value = -271.1539199999997
$ws.Range("N137").Value = -7696.0713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4457.2754
$ws.Range("I32").Value = 3353.1475
$ws.Range("J32").Value = 12876.25
$ws.Range("K32").Value = 3353.1475
$ws.Range("L32").Value = 12876.25
$ws.Range("M32").Value = -3066.1475
$ws.Range("N32").Value = -13450.25

$ws.Range("H45").Value = 17902
$ws.Range("I45").Value = 17902
$ws.Range("K45").Value = 17902
$ws.Range("M45").Value = -17525

$ws.Range("H61").Value = 3266.1914
$ws.Range("I61").Value = 3506.9285
$ws.Range("J61").Value = 1244
$ws.Range("K61").Value = 3506.9285
$ws.Range("L61").Value = 1244
$ws.Range("M61").Value = -3294.9285
$ws.Range("N61").Value = -1668

$ws.Range("H74").Value = 1155.5625
$ws.Range("I74").Value = 1090.3478
$ws.Range("J74").Value = 1322.2222
$ws.Range("K74").Value = 1090.3478
$ws.Range("L74").Value = 1322.2222
$ws.Range("M74").Value = -216.3478
$ws.Range("N74").Value = -3070.2222

$ws.Range("H77").Value = 1155.5625
$ws.Range("I77").Value = 1090.3478
$ws.Range("J77").Value = 1322.2222
$ws.Range("K77").Value = 5451.739
$ws.Range("L77").Value = 6611.111
$ws.Range("M77").Value = -1083.739
$ws.Range("N77").Value = -15347.111

$ws.Range("H132").Value = 2452.4473
$ws.Range("I132").Value = 1241.6818
$ws.Range("J132").Value = 4117.25
$ws.Range("K132").Value = 3725.0454
$ws.Range("L132").Value = 12351.75
$ws.Range("M132").Value = -1195.0454
$ws.Range("N132").Value = -17411.75

$ws.Range("H136").Value = 3266.1914
$ws.Range("I136").Value = 3506.9285
$ws.Range("J136").Value = 1244
$ws.Range("K136").Value = 10520.7855
$ws.Range("L136").Value = 3732
$ws.Range("M136").Value = -7970.7855
$ws.Range("N136").Value = -8832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 55556456
$ws.Range("I99").Value = 83333870
$ws.Range("J99").Value = 1624.6666
$ws.Range("K99").Value = 83333870
$ws.Range("L99").Value = 1624.6666
$ws.Range("M99").Value = -83332372
$ws.Range("N99").Value = -4620.6666

$ws.Range("H107").Value = 1019.375
$ws.Range("I107").Value = 1060.6666
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1060.6666
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 859.3334
$ws.Range("N107").Value = -4240

$ws.Range("H134").Value = 4123.439
$ws.Range("I134").Value = 4426.091
$ws.Range("K134").Value = 13278.273
$ws.Range("M134").Value = -10743.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3476.0789
$ws.Range("I31").Value = 1676.5625
$ws.Range("J31").Value = 4784.8184
$ws.Range("K31").Value = 1676.5625
$ws.Range("L31").Value = 4784.8184
$ws.Range("M31").Value = -1381.5625
$ws.Range("N31").Value = -5374.8184

$ws.Range("H34").Value = 3476.0789
$ws.Range("I34").Value = 1676.5625
$ws.Range("J34").Value = 4784.8184
$ws.Range("K34").Value = 1676.5625
$ws.Range("L34").Value = 4784.8184
$ws.Range("M34").Value = -1474.5625
$ws.Range("N34").Value = -5188.8184

$ws.Range("H58").Value = 1119.804
$ws.Range("I58").Value = 763.51514
$ws.Range("J58").Value = 1773
$ws.Range("K58").Value = 763.51514
$ws.Range("L58").Value = 1773
$ws.Range("M58").Value = -560.51514
$ws.Range("N58").Value = -2179

$ws.Range("H99").Value = 4811.8335
$ws.Range("I99").Value = 3333.2666
$ws.Range("J99").Value = 12204.667
$ws.Range("K99").Value = 3333.2666
$ws.Range("L99").Value = 12204.667
$ws.Range("M99").Value = -1835.2666
$ws.Range("N99").Value = -15200.667

$ws.Range("H107").Value = 262.64706
$ws.Range("I107").Value = 188.63637
$ws.Range("J107").Value = 298.0435
$ws.Range("K107").Value = 188.63637
$ws.Range("L107").Value = 298.0435
$ws.Range("M107").Value = 1731.36363
$ws.Range("N107").Value = -4138.0435

$ws.Range("H126").Value = 4811.8335
$ws.Range("I126").Value = 3333.2666
$ws.Range("J126").Value = 12204.667
$ws.Range("K126").Value = 9999.799800000001
$ws.Range("L126").Value = 36614.001
$ws.Range("M126").Value = -7529.799800000001
$ws.Range("N126").Value = -41554.001

$ws.Range("H132").Value = 2067.9556
$ws.Range("I132").Value = 1854.2858
$ws.Range("J132").Value = 2815.8
$ws.Range("K132").Value = 5562.857400000001
$ws.Range("L132").Value = 8447.400000000001
$ws.Range("M132").Value = -3032.857400000001
$ws.Range("N132").Value = -13507.4

$ws.Range("H134").Value = 1971.0667
$ws.Range("I134").Value = 2133.4285
$ws.Range("K134").Value = 6400.2855
$ws.Range("M134").Value = -3865.2855

$ws.Range("H136").Value = 1119.804
$ws.Range("I136").Value = 763.51514
$ws.Range("J136").Value = 1773
$ws.Range("K136").Value = 2290.54542
$ws.Range("L136").Value = 5319
$ws.Range("M136").Value = 259.4545800000001
$ws.Range("N136").Value = -10419

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1023.5455
$ws.Range("I102").Value = 965.9
$ws.Range("K102").Value = 965.9
$ws.Range("M102").Value = 656.1

$ws.Range("H107").Value = 2207.3076
$ws.Range("I107").Value = 669.8570999999999
$ws.Range("J107").Value = 4001
$ws.Range("K107").Value = 669.8570999999999
$ws.Range("L107").Value = 4001
$ws.Range("M107").Value = 1250.1429
$ws.Range("N107").Value = -7841

$ws.Range("H126").Value = 6283.7393
$ws.Range("I126").Value = 8907.929
$ws.Range("J126").Value = 2201.6667
$ws.Range("K126").Value = 26723.787
$ws.Range("L126").Value = 6605.000100000001
$ws.Range("M126").Value = -24253.787
$ws.Range("N126").Value = -11545.0001

$ws.Range("H132").Value = 2799.7551
$ws.Range("I132").Value = 2801.2258
$ws.Range("J132").Value = 2797.2222
$ws.Range("K132").Value = 8403.6774
$ws.Range("L132").Value = 8391.6666
$ws.Range("M132").Value = -5873.6774
$ws.Range("N132").Value = -13451.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1656.3182
$ws.Range("I7").Value = 1375.8182
$ws.Range("J7").Value = 1936.8182
$ws.Range("K7").Value = 1375.8182
$ws.Range("L7").Value = 1936.8182
$ws.Range("M7").Value = -1263.8182
$ws.Range("N7").Value = -2160.8182

$ws.Range("H40").Value = 83335720
$ws.Range("I40").Value = 166668980
$ws.Range("J40").Value = 2451.6667
$ws.Range("K40").Value = 166668980
$ws.Range("L40").Value = 2451.6667
$ws.Range("M40").Value = -166668844
$ws.Range("N40").Value = -2723.6667

$ws.Range("H126").Value = 1656.3182
$ws.Range("I126").Value = 1375.8182
$ws.Range("J126").Value = 1936.8182
$ws.Range("K126").Value = 4127.4546
$ws.Range("L126").Value = 5810.4546
$ws.Range("M126").Value = -1657.4546
$ws.Range("N126").Value = -10750.4546

$ws.Range("H132").Value = 6281608.5
$ws.Range("I132").Value = 9755450
$ws.Range("J132").Value = 1972.3462
$ws.Range("K132").Value = 29266350
$ws.Range("L132").Value = 5917.0386
$ws.Range("M132").Value = -29263820
$ws.Range("N132").Value = -10977.0386

$ws.Range("H136").Value = 6319.237
$ws.Range("I136").Value = 6615.2607
$ws.Range("J136").Value = 5865.3335
$ws.Range("K136").Value = 19845.7821
$ws.Range("L136").Value = 17596.0005
$ws.Range("M136").Value = -17295.7821
$ws.Range("N136").Value = -22696.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 38461910
$ws.Range("I107").Value = 52631916
$ws.Range("J107").Value = 472.14285
$ws.Range("K107").Value = 157895748
$ws.Range("L107").Value = 1416.42855
$ws.Range("M107").Value = -157893828
$ws.Range("N107").Value = -5256.428550000001

$ws.Range("H126").Value = 1361.6842
$ws.Range("I126").Value = 1159.5454
$ws.Range("J126").Value = 1639.625
$ws.Range("K126").Value = 3478.6362
$ws.Range("L126").Value = 4918.875
$ws.Range("M126").Value = -1008.6362
$ws.Range("N126").Value = -9858.875

$ws.Range("H132").Value = 1309.6938
$ws.Range("I132").Value = 871.3823
$ws.Range("J132").Value = 2303.2
$ws.Range("K132").Value = 2614.1469
$ws.Range("L132").Value = 6909.599999999999
$ws.Range("M132").Value = -84.14689999999973
$ws.Range("N132").Value = -11969.6

$ws.Range("H136").Value = 1113.1936
$ws.Range("I136").Value = 633.75
$ws.Range("J136").Value = 1984.909
$ws.Range("K136").Value = 1901.25
$ws.Range("L136").Value = 5954.727000000001
$ws.Range("M136").Value = 648.75
$ws.Range("N136").Value = -11054.727
